$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.418.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.683.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.85%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.369"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.157.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.253.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("E16").Value = "  -3.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.685.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.65%  "

$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.507"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0860"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.02%  "

$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("E31").Value = "  -4.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.52%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -2.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "343.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.47%  "

$ws.Range("E39").Value = "  -5.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.76%  "

$ws.Range("E42").Value = "  -6.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.619"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0563"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.37%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.64%  "

$ws.Range("E51").Value = "  -4.68%  "
